# Auto-applied numeric updates to H:N profit columns across 8 sheets
# per the scheduled-runner data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 52
$ws.Range("H52").Value = 402.8
$ws.Range("J52").Value = 402.8
$ws.Range("L52").Value = 1208.4
$ws.Range("N52").Value = -1528.4

# Row 70
$ws.Range("H70").Value = 8750
$ws.Range("I70").Value = 5000
$ws.Range("K70").Value = 15000
$ws.Range("M70").Value = -14730

# Row 73
$ws.Range("H73").Value = 8750
$ws.Range("I73").Value = 5000
$ws.Range("K73").Value = 15000
$ws.Range("M73").Value = -14064

# Row 86
$ws.Range("H86").Value = 5000
$ws.Range("I86").Value = 5000
$ws.Range("K86").Value = 5000
$ws.Range("M86").Value = -3877

# Row 88
$ws.Range("H88").Value = 1440.8
$ws.Range("I88").Value = 902.5
$ws.Range("J88").Value = 1799.6666
$ws.Range("K88").Value = 902.5
$ws.Range("L88").Value = 1799.6666
$ws.Range("M88").Value = -496.5
$ws.Range("N88").Value = -2611.6666

# Row 89
$ws.Range("H89").Value = 5000
$ws.Range("I89").Value = 5000
$ws.Range("K89").Value = 25000
$ws.Range("M89").Value = -19384

# Row 91
$ws.Range("H91").Value = 1440.8
$ws.Range("I91").Value = 902.5
$ws.Range("J91").Value = 1799.6666
$ws.Range("K91").Value = 902.5
$ws.Range("L91").Value = 1799.6666
$ws.Range("M91").Value = 501.5
$ws.Range("N91").Value = -4607.6666

# Row 94
$ws.Range("H94").Value = 7210.5
$ws.Range("I94").Value = 7210.5
$ws.Range("K94").Value = 7210.5
$ws.Range("M94").Value = -6759.5

# Row 132
$ws.Range("H132").Value = 15783.467
$ws.Range("I132").Value = 20650.273
$ws.Range("J132").Value = 2399.75
$ws.Range("K132").Value = 61950.819
$ws.Range("L132").Value = 7199.25
$ws.Range("M132").Value = -59420.819
$ws.Range("N132").Value = -12259.25

$ws = $wb.Worksheets.Item("ARM")
# Row 63
$ws.Range("H63").Value = 3415.6924
$ws.Range("J63").Value = 4086
$ws.Range("L63").Value = 4086
$ws.Range("N63").Value = -5458

# Row 66
$ws.Range("H66").Value = 3415.6924
$ws.Range("J66").Value = 4086
$ws.Range("L66").Value = 20430
$ws.Range("N66").Value = -27294

# Row 88
$ws.Range("H88").Value = 2478.5715
$ws.Range("I88").Value = 1987.5
$ws.Range("J88").Value = 2675
$ws.Range("K88").Value = 1987.5
$ws.Range("L88").Value = 2675
$ws.Range("M88").Value = -1581.5
$ws.Range("N88").Value = -3487

# Row 91
$ws.Range("H91").Value = 2478.5715
$ws.Range("I91").Value = 1987.5
$ws.Range("J91").Value = 2675
$ws.Range("K91").Value = 1987.5
$ws.Range("L91").Value = 2675
$ws.Range("M91").Value = -583.5
$ws.Range("N91").Value = -5483

# Row 132
$ws.Range("H132").Value = 3961.2856
$ws.Range("J132").Value = 7316.3335
$ws.Range("L132").Value = 21949.0005
$ws.Range("N132").Value = -27009.0005

$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("H94").Value = 1700
$ws.Range("I94").Value = 1700
$ws.Range("K94").Value = 1700
$ws.Range("M94").Value = -1249

$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Range("H7").Value = 96
$ws.Range("I7").Value = 65.333336
$ws.Range("K7").Value = 65.333336
$ws.Range("M7").Value = 47.666664

# Row 22
$ws.Range("H22").Value = 1462.5834
$ws.Range("I22").Value = 281.375
$ws.Range("J22").Value = 3825
$ws.Range("K22").Value = 281.375
$ws.Range("L22").Value = 3825
$ws.Range("M22").Value = 68.625
$ws.Range("N22").Value = -4525

# Row 105
$ws.Range("H105").Value = 1754.6428
$ws.Range("I105").Value = 1754.6428
$ws.Range("K105").Value = 1754.6428
$ws.Range("M105").Value = -7.642800000000079

$ws = $wb.Worksheets.Item("CUL")
# Row 107
$ws.Range("H107").Value = 932.8889
$ws.Range("J107").Value = 874.2143
$ws.Range("L107").Value = 2622.6429
$ws.Range("N107").Value = -6462.6429

$ws = $wb.Worksheets.Item("GSM")
# Row 40
$ws.Range("H40").Value = 9999
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 9999
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 9999
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = -10301

# Row 44
$ws.Range("H44").Value = 4000
$ws.Range("J44").Value = 4000
$ws.Range("L44").Value = 4000
$ws.Range("N44").Value = -5192

# Row 80
$ws.Range("H80").Value = 2863.1538
$ws.Range("I80").Value = 3065.818
$ws.Range("J80").Value = 1748.5
$ws.Range("K80").Value = 3065.818
$ws.Range("L80").Value = 1748.5
$ws.Range("M80").Value = -2067.818
$ws.Range("N80").Value = -3744.5

# Row 83
$ws.Range("H83").Value = 2863.1538
$ws.Range("I83").Value = 3065.818
$ws.Range("J83").Value = 1748.5
$ws.Range("K83").Value = 15329.09
$ws.Range("L83").Value = 8742.5
$ws.Range("M83").Value = -10337.09
$ws.Range("N83").Value = -18726.5

# Row 122
$ws.Range("H122").Value = 195551.11
$ws.Range("I122").Value = 296860.2
$ws.Range("J122").Value = 4189.5557
$ws.Range("K122").Value = 890580.6000000001
$ws.Range("L122").Value = 12568.6671
$ws.Range("M122").Value = -888130.6000000001
$ws.Range("N122").Value = -17468.6671

# Row 136
$ws.Range("H136").Value = 30348.715
$ws.Range("J136").Value = 30348.715
$ws.Range("L136").Value = 91046.145
$ws.Range("N136").Value = -96146.145

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 2966.3333
$ws.Range("I22").Value = 1950
$ws.Range("J22").Value = 4999
$ws.Range("K22").Value = 1950
$ws.Range("L22").Value = 4999
$ws.Range("M22").Value = -1655
$ws.Range("N22").Value = -5589

# Row 27
$ws.Range("H27").Value = 2966.3333
$ws.Range("I27").Value = 1950
$ws.Range("J27").Value = 4999
$ws.Range("K27").Value = 1950
$ws.Range("L27").Value = 4999
$ws.Range("M27").Value = -1843
$ws.Range("N27").Value = -5213

# Row 46
$ws.Range("H46").Value = 5622.75
$ws.Range("I46").Value = 4096.8
$ws.Range("J46").Value = 6316.364
$ws.Range("K46").Value = 4096.8
$ws.Range("L46").Value = 6316.364
$ws.Range("M46").Value = -3908.8
$ws.Range("N46").Value = -6692.364

# Row 122
$ws.Range("H122").Value = 3730
$ws.Range("I122").Value = 3391.6667
$ws.Range("J122").Value = 4745
$ws.Range("K122").Value = 10175.0001
$ws.Range("L122").Value = 14235
$ws.Range("M122").Value = -7725.000100000001
$ws.Range("N122").Value = -19135

# Row 130
$ws.Range("H130").Value = 30000
$ws.Range("J130").Value = 30000
$ws.Range("L130").Value = 30000
$ws.Range("N130").Value = -40040

# Row 136
$ws.Range("H136").Value = 4278.2
$ws.Range("I136").Value = 4278.2
$ws.Range("K136").Value = 12834.6
$ws.Range("M136").Value = -10284.6

$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 0
$ws.Range("I81").Value = 0
$ws.Range("K81").Value = 0
$ws.Range("M81").ClearContents()

# Row 84
$ws.Range("H84").Value = 0
$ws.Range("I84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("M84").ClearContents()

# Row 122
$ws.Range("H122").Value = 1226.6111
$ws.Range("I122").Value = 1226.6111
$ws.Range("K122").Value = 3679.8333
$ws.Range("M122").Value = -1229.8333

# Row 136
$ws.Range("H136").Value = 2542.889
$ws.Range("I136").Value = 1524.5358
$ws.Range("K136").Value = 4573.607400000001
$ws.Range("M136").Value = -2023.607400000001
